$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new column "S" (year 2022) to the table, mirroring the formatting
# of the existing "R" column (year 2021 / value).

# Header cell S4: copy formatting from R4 (bold, right/center aligned,
# bordered) and set the new year value.
$ws.Range("R4").Copy()
$ws.Range("S4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("S4").Value = 2022

# Data cell S5: copy formatting from R5, then tweak font (regular weight,
# size 9) and number format (one decimal place) to match the new style,
# and set the new data value.
$ws.Range("R5").Copy()
$ws.Range("S5").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("S5").Value = 42
$ws.Range("S5").Font.Bold = $false
$ws.Range("S5").Font.Size = 9
$ws.Range("S5").NumberFormat = "0.0"

$excel.CutCopyMode = 0

# Update the active selection to reflect where the user ended up editing.
[void]$ws.Range("U4").Select()
